# Rewrite login test case: replace the "Tittle" column header with
# "Expected_tittle", update the selected cell, and widen column C so the
# new header fits better.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1 header text changed from "Tittle" to "Expected_tittle"
$ws.Range("C1").Value = "Expected_tittle"

# Selection moved from E5 to D10
$ws.Range("D10").Select() | Out-Null

# Column C (previously grouped with column B at width 12.796875) now has
# its own, wider column width so the longer header text fits.
$ws.Columns("C").ColumnWidth = 14.8
